$d = $word.ActiveDocument

# 1) Fix the text of the first (existing) paragraph: drop the "(хотя бы ШИМ-каналы)" qualifier
#    and the trailing "." run, merging into a single run with the trimmed text.
$found = $d.Content.Find.Execute(
    "Выровнять длины всех дорожек (хотя бы ШИМ-каналы).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Выровнять длины всех дорожек", 2
)
if (-not $found) {
    throw "Could not find the first bullet's original text to rewrite it."
}

# 2) Append the new bullet paragraphs (same list style/numbering as paragraph 1) by
#    inserting a WordOpenXML fragment at the end of the document content. Using raw
#    OOXML here guarantees exact run-splitting / rPr (w:lang, w:vertAlign) fidelity.
$ins = $d.Content
$ins.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Питание 3.3 В вместо 5 В</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Подумать над решением проблемы с индуктивностью – дают большое падение напряжения</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Соблюсти правила по </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>ESL</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>на каждой индуктивности, чтобы не было резонанса</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Создать правила (хотя бы для </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>I</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:vertAlign w:val="superscript"/>
              </w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>C</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>линии)</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Переделать полигон с землёй (полигон под кварцем сделать под размер корпуса)</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="a3"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="3"/>
              </w:numPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$res = $ins.InsertXML($xml)
if (-not $res) {
    throw "InsertXML failed to append the new Rev.2 bullet paragraphs."
}

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
